# Varios/traducciones.xlsx — add 5 new translated messages (rows 168-172)
# to the "Controles" sheet: com.td.horario.no.disponible, com.td.curso.codigo.existe,
# com.td.curso.tiene.alumnos, com.td.repetidos, com.td.familia.existe.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Controles")

# Codigo / Español / Ingles for the new messages, in the order they were appended.
$newRows = @(
    @{ Row = 168; Codigo = "com.td.horario.no.disponible"; Es = "Horario no disponible";            En = "Schedule not available" },
    @{ Row = 169; Codigo = "com.td.curso.codigo.existe";    Es = "Código existente";                 En = "The code already exists" },
    @{ Row = 170; Codigo = "com.td.curso.tiene.alumnos";    Es = "Tiene alumnos";                    En = "It has students" },
    @{ Row = 171; Codigo = "com.td.repetidos";              Es = "El mail o el D.N.I. ya existe";    En = "The mail or the ID already exist" },
    @{ Row = 172; Codigo = "com.td.familia.existe";         Es = "El nombre de la familia existe";   En = "The family name already exists" }
)

foreach ($item in $newRows) {
    $r = $item.Row

    $ws.Cells.Item($r, 1).Value = $item.Codigo
    $ws.Cells.Item($r, 2).Value = $item.Es
    $ws.Cells.Item($r, 3).Value = $item.En

    # Same CLEAN(CONCAT(...)) "insert into MENSAJE..." builder formulas used by every
    # other row in the sheet (D = Spanish insert statement, E = English insert statement).
    $ws.Cells.Item($r, 4).Formula = '=CLEAN(CONCAT("insert into MENSAJE(MSJ_CODIGO,MSJ_IDIOMA_ID,MSJ_TEXTO)
values (''",$A' + $r + ',"'',( select idioma.IDI_ID from IDIOMA where IDI_CODIGO = ''",$D$1,"''),''",$B' + $r + ',"'')"))'

    $ws.Cells.Item($r, 5).Formula = '=CLEAN(CONCAT("insert into MENSAJE(MSJ_CODIGO,MSJ_IDIOMA_ID,MSJ_TEXTO)
values (''",$A' + $r + ',"'',( select idioma.IDI_ID from IDIOMA where IDI_CODIGO = ''",$E$1,"''),''",$C' + $r + ',"'')"))'

    # Entering a formula with an embedded line break auto-grows the row; put it back to
    # the sheet's normal (non-custom) height like every other data row.
    $ws.Rows.Item($r).AutoFit()
}

# Move the on-screen selection the way the author's session ended up (cursor on A42).
$ws.Activate() | Out-Null
$ws.Range("A42").Select() | Out-Null
